$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.02679002569074385
$ws.Range("D2").Value = 0.2318299200193223
$ws.Range("E2").Value = 0.162938907540493
$ws.Range("F2").Value = 1.004500738417136
$ws.Range("G2").Value = 0.4721807277463057
$ws.Range("H2").Value = 0.6335810609531052
$ws.Range("I2").Value = 0.4738580892836488
$ws.Range("J2").Value = 0.1586882765817705
$ws.Range("K2").Value = 1.250536817984766
$ws.Range("M2").Value = 0.3901111461669942
$ws.Range("O2").Value = 2.170523178034159
$ws.Range("C3").Value = 0.02357037927475858
$ws.Range("D3").Value = 0.228962982107916
$ws.Range("E3").Value = 0.1632560125944806
$ws.Range("F3").Value = 1.013960101197725
$ws.Range("G3").Value = 0.4788143987502309
$ws.Range("H3").Value = 0.6411106823593258
$ws.Range("I3").Value = 0.4785490190033599
$ws.Range("J3").Value = 0.1605026073238101
$ws.Range("K3").Value = 1.096027920285394
$ws.Range("M3").Value = 0.3571328374326157
$ws.Range("O3").Value = 2.200075201232522
$ws.Range("C4").Value = 0.02158368407354061
$ws.Range("D4").Value = 0.2272829819026896
$ws.Range("E4").Value = 0.1635314992623407
$ws.Range("F4").Value = 1.020462444954624
$ws.Range("G4").Value = 0.4833321012252227
$ws.Range("H4").Value = 0.6460864198051226
$ws.Range("I4").Value = 0.4818082336998515
$ws.Range("J4").Value = 0.1617136677624025
$ws.Range("K4").Value = 1.000794740260488
$ws.Range("M4").Value = 0.3368782997884523
$ws.Range("O4").Value = 2.219891462497358
$ws.Range("C5").Value = 0.02077166140531261
$ws.Range("D5").Value = 0.2266186588779391
$ws.Range("E5").Value = 0.1636640963411544
$ws.Range("F5").Value = 1.023286658724381
$ws.Range("G5").Value = 0.4852846602912066
$ws.Range("H5").Value = 0.6482026930336389
$ws.Range("I5").Value = 0.4832315644271006
$ws.Range("J5").Value = 0.162231576708848
$ws.Range("K5").Value = 0.9618977914338132
$ws.Range("M5").Value = 0.3286235808157087
$ws.Range("O5").Value = 2.2283864212189
$ws.Range("C6").Value = 0.02063668020738163
$ws.Range("D6").Value = 0.2265095768778167
$ws.Range("E6").Value = 0.163687342713466
$ws.Range("F6").Value = 1.023766150422261
$ws.Range("G6").Value = 0.4856156120414994
$ws.Range("H6").Value = 0.6485594503021019
$ws.Range("I6").Value = 0.4834736521824823
$ws.Range("J6").Value = 0.1623190479273333
$ws.Range("K6").Value = 0.9554337139689721
$ws.Range("M6").Value = 0.3272528588542656
$ws.Range("O6").Value = 2.229822333130556
$ws.Range("C7").Value = 0.02157274262954445
$ws.Range("D7").Value = 0.2272739403257162
$ws.Range("E7").Value = 0.1635332051539109
$ws.Range("F7").Value = 1.020499827085025
$ws.Range("G7").Value = 0.4833579827423833
$ws.Range("H7").Value = 0.6461146018204289
$ws.Range("I7").Value = 0.4818270440380878
$ws.Range("J7").Value = 0.1617205537187623
$ws.Range("K7").Value = 1.000270517504816
$ws.Range("M7").Value = 0.3367669763140881
$ws.Range("O7").Value = 2.220004330001487
$ws.Range("C8").Value = 0.0256819531550434
$ws.Range("D8").Value = 0.2308247872898335
$ws.Range("E8").Value = 0.1630314908689314
$ws.Range("F8").Value = 1.007618161075705
$ws.Range("G8").Value = 0.4743755965954435
$ws.Range("H8").Value = 0.63610409505538
$ws.Range("I8").Value = 0.4753967958012559
$ws.Range("J8").Value = 0.1592937065080875
$ws.Range("K8").Value = 1.197339402928662
$ws.Range("M8").Value = 0.3787417667015376
$ws.Range("O8").Value = 2.180365520408841
$ws.Range("C9").Value = 0.03366076766420179
$ws.Range("D9").Value = 0.238421555118137
$ws.Range("E9").Value = 0.1626878559824014
$ws.Range("F9").Value = 0.9878716617470857
$ws.Range("G9").Value = 0.4602996963741646
$ws.Range("H9").Value = 0.6192715918664433
$ws.Range("I9").Value = 0.4657995945620712
$ws.Range("J9").Value = 0.1553053321431932
$ws.Range("K9").Value = 1.580793540683544
$ws.Range("M9").Value = 0.4609862174437609
$ws.Range("O9").Value = 2.115919426327793
$ws.Range("C10").Value = 0.03947307542422607
$ws.Range("D10").Value = 0.2443850109013965
$ws.Range("E10").Value = 0.162824899877819
$ws.Range("F10").Value = 0.976734169667381
$ws.Range("G10").Value = 0.4521302698120166
$ws.Range("H10").Value = 0.6086117563888322
$ws.Range("I10").Value = 0.4605930244536225
$ws.Range("J10").Value = 0.1528456993283243
$ws.Range("K10").Value = 1.86057213172711
$ws.Range("M10").Value = 0.521345545996553
$ws.Range("O10").Value = 2.076704727688849
$ws.Range("C11").Value = 0.04210619249512604
$ws.Range("D11").Value = 0.2471800908597288
$ws.Range("E11").Value = 0.1629716718827616
$ws.Range("F11").Value = 0.9724011141706512
$ws.Range("G11").Value = 0.448888640381405
$ws.Range("H11").Value = 0.6041332457214281
$ws.Range("I11").Value = 0.4586266590234516
$ws.Range("J11").Value = 0.1518291155635865
$ws.Range("K11").Value = 1.987405336701613
$ws.Range("M11").Value = 0.5487856039385832
$ws.Range("O11").Value = 2.060638543800309
$ws.Range("C12").Value = 0.04310167963352285
$ws.Range("D12").Value = 0.2482502607645642
$ws.Range("E12").Value = 0.163039374896357
$ws.Range("F12").Value = 0.9708659119991339
$ws.Range("G12").Value = 0.4477296380635636
$ws.Range("H12").Value = 0.6024906868957274
$ws.Range("I12").Value = 0.4579400081889347
$ws.Range("J12").Value = 0.1514588901283496
$ws.Range("K12").Value = 2.035368082544608
$ws.Range("M12").Value = 0.5591733717291021
$ws.Range("O12").Value = 2.054810240842173
$ws.Range("C13").Value = 0.04288735631050145
$ws.Range("D13").Value = 0.248019260208082
$ws.Range("E13").Value = 0.1630242549152783
$ws.Range("F13").Value = 0.9711918446676648
$ws.Range("G13").Value = 0.4479761976196954
$ws.Range("H13").Value = 0.6028420675223671
$ws.Range("I13").Value = 0.4580853102746758
$ws.Range("J13").Value = 0.1515379692396515
$ws.Range("K13").Value = 2.025041441051258
$ws.Range("M13").Value = 0.5569363319572886
$ws.Range("O13").Value = 2.056054090828837
$ws.Range("C14").Value = 0.042188124492057
$ws.Range("D14").Value = 0.247267899774485
$ws.Range("E14").Value = 0.1629769989739955
$ws.Range("F14").Value = 0.9722726942146735
$ws.Range("G14").Value = 0.4487919135123022
$ws.Range("H14").Value = 0.6039970418031331
$ws.Range("I14").Value = 0.4585690052435112
$ws.Range("J14").Value = 0.1517983615138796
$ws.Range("K14").Value = 1.991352605456939
$ws.Range("M14").Value = 0.5496402793663577
$ws.Range("O14").Value = 2.060153919031194
$ws.Range("C15").Value = 0.04175961286885865
$ws.Range("D15").Value = 0.2468091950217257
$ws.Range("E15").Value = 0.1629496316865584
$ws.Range("F15").Value = 0.9729485069776871
$ws.Range("G15").Value = 0.4493004960735476
$ws.Range("H15").Value = 0.6047114468850339
$ws.Range("I15").Value = 0.4588728358929188
$ws.Range("J15").Value = 0.1519597785314168
$ws.Range("K15").Value = 1.970708505240566
$ws.Range("M15").Value = 0.5451708025481139
$ws.Range("O15").Value = 2.0626984907977
$ws.Range("C16").Value = 0.03930077119944997
$ws.Range("D16").Value = 0.2442039931078739
$ws.Range("E16").Value = 0.1628170054538245
$ws.Range("F16").Value = 0.977032117839812
$ws.Range("G16").Value = 0.4523516931886604
$ws.Range("H16").Value = 0.608911901299777
$ws.Range("I16").Value = 0.4607296352192591
$ws.Range("J16").Value = 0.1529141960243194
$ws.Range("K16").Value = 1.852274179597543
$ws.Range("M16").Value = 0.5195518625504718
$ws.Range("O16").Value = 2.077790422420961
$ws.Range("C17").Value = 0.0377895195409792
$ws.Range("D17").Value = 0.2426267892749649
$ws.Range("E17").Value = 0.162757255834574
$ws.Range("F17").Value = 0.9797252506570473
$ws.Range("G17").Value = 0.45434527686065
$ws.Range("H17").Value = 0.6115837249851666
$ws.Range("I17").Value = 0.4619718189259032
$ws.Range("J17").Value = 0.153525919922668
$ws.Range("K17").Value = 1.77950382438479
$ws.Range("M17").Value = 0.5038304809518053
$ws.Range("O17").Value = 2.087503376461427
$ws.Range("C18").Value = 0.03691926139219959
$ws.Range("D18").Value = 0.2417273731553138
$ws.Range("E18").Value = 0.1627308374296099
$ws.Range("F18").Value = 0.9813432951427004
$ws.Range("G18").Value = 0.4555365967101608
$ws.Range("H18").Value = 0.6131553765818367
$ws.Range("I18").Value = 0.462724131987013
$ws.Range("J18").Value = 0.1538873965696617
$ws.Range("K18").Value = 1.737607054868874
$ws.Range("M18").Value = 0.4947863297692408
$ws.Range("O18").Value = 2.093256832218614
$ws.Range("C19").Value = 0.03662443207345234
$ws.Range("D19").Value = 0.2414241803112276
$ws.Range("E19").Value = 0.162723258213088
$ws.Range("F19").Value = 0.9819029874307788
$ws.Range("G19").Value = 0.4559476197349142
$ws.Range("H19").Value = 0.613693501564903
$ws.Range("I19").Value = 0.4629853472453789
$ws.Range("J19").Value = 0.1540114395734804
$ws.Range("K19").Value = 1.723414547939228
$ws.Range("M19").Value = 0.4917238755191136
$ws.Range("O19").Value = 2.095233484174358
$ws.Range("C20").Value = 0.03795050139204648
$ws.Range("D20").Value = 0.2427938839448274
$ws.Range("E20").Value = 0.1627627938210416
$ws.Range("F20").Value = 0.9794314169782083
$ws.Range("G20").Value = 0.4541284320090426
$ws.Range("H20").Value = 0.6112956935233456
$ws.Range("I20").Value = 0.4618356688171694
$ws.Range("J20").Value = 0.1534598041346591
$ws.Range("K20").Value = 1.787254637615661
$ws.Range("M20").Value = 0.5055042220417363
$ws.Range("O20").Value = 2.086452145517967
$ws.Range("C21").Value = 0.0423935500747632
$ws.Range("D21").Value = 0.2474882748683314
$ws.Range("E21").Value = 0.1629905502955431
$ws.Range("F21").Value = 0.9719523543909574
$ws.Range("G21").Value = 0.4485504558625379
$ws.Range("H21").Value = 0.6036563494008647
$ws.Range("I21").Value = 0.4584253579069184
$ws.Range("J21").Value = 0.1517214779894012
$ws.Range("K21").Value = 2.001249654824392
$ws.Range("M21").Value = 0.5517833984405911
$ws.Range("O21").Value = 2.058942757399976
$ws.Range("C22").Value = 0.04528789230201369
$ws.Range("D22").Value = 0.2506246793177667
$ws.Range("E22").Value = 0.163210060090119
$ws.Range("F22").Value = 0.9676801326532427
$ws.Range("G22").Value = 0.4453045107680893
$ws.Range("H22").Value = 0.5989746185869436
$ws.Range("I22").Value = 0.4565344762192396
$ws.Range("J22").Value = 0.1506712644024617
$ws.Range("K22").Value = 2.140720301982356
$ws.Range("M22").Value = 0.5820106471873316
$ws.Range("O22").Value = 2.042454001770821
$ws.Range("C23").Value = 0.0437440067165511
$ws.Range("D23").Value = 0.2489444993508982
$ws.Range("E23").Value = 0.163086444126499
$ws.Range("F23").Value = 0.9699039006028869
$ws.Range("G23").Value = 0.4470002808715066
$ws.Range("H23").Value = 0.6014448716034977
$ws.Range("I23").Value = 0.4575127055973169
$ws.Range("J23").Value = 0.1512239183044564
$ws.Range("K23").Value = 2.066318643445015
$ws.Range("M23").Value = 0.5658797327065912
$ws.Range("O23").Value = 2.051117779431749
$ws.Range("C24").Value = 0.03787772591526561
$ws.Range("D24").Value = 0.2427183175767169
$ws.Range("E24").Value = 0.162760265387444
$ws.Range("F24").Value = 0.9795640420035028
$ws.Range("G24").Value = 0.4542263268484135
$ws.Range("H24").Value = 0.6114258017068579
$ws.Range("I24").Value = 0.4618971034107204
$ws.Range("J24").Value = 0.1534896645938026
$ws.Range("K24").Value = 1.783750682001539
$ws.Range("M24").Value = 0.5047475414626774
$ws.Range("O24").Value = 2.086926879567216
$ws.Range("C25").Value = 0.0315109216842302
$ws.Range("D25").Value = 0.2362990283982924
$ws.Range("E25").Value = 0.1627123588362167
$ws.Range("F25").Value = 0.9926223695293288
$ws.Range("G25").Value = 0.4637272238860888
$ws.Range("H25").Value = 0.6235255498984031
$ws.Range("I25").Value = 0.4680725519228268
$ws.Range("J25").Value = 0.1563017041857329
$ws.Range("K25").Value = 1.477392217004478
$ws.Range("M25").Value = 0.438746714777821
$ws.Range("O25").Value = 2.13192790136506
